{"js": "// Fix the missing period at the end of Anthony's \"What could be improved?\"\n// paragraph and move the \"_GoBack\" bookmark from the end of the document to\n// right after the newly-added period (this mirrors what Word itself does:\n// \"_GoBack\" always tracks the location of the most recent edit).\n\n// 1) Remove the existing \"_GoBack\" bookmark (currently sitting at the very\n//    end of the document, after the last paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the paragraph that ends in \"...which is not ideal\" (the\n//    \"What could be improved?\" answer) and append the missing \".\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"which is not ideal\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nconst endRange = target.getRange(\"End\");\n\n// Use insertOoxml so the new run carries the exact same run formatting\n// (w:lang w:val=\"en-AU\") as the rest of the paragraph's runs.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:rPr><w:lang w:val=\"en-AU\"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nendRange.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark right after the period we just added.\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet newTarget = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"which is not ideal\") !== -1) {\n    newTarget = p;\n    break;\n  }\n}\n\nnewTarget.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Fix the missing period at the end of Anthony's \"What could be improved?\"\n# paragraph and move the \"_GoBack\" bookmark from the end of the document to\n# right after the newly-added period (mirrors Word's own behaviour of\n# keeping \"_GoBack\" anchored at the location of the most recent edit).\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (currently sitting at the very\n#    end of the document, after the last paragraph). \"_GoBack\" is a hidden\n#    bookmark, but it is still addressable by name directly.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Locate the paragraph that ends in \"...which is not ideal\" (the\n#    \"What could be improved?\" answer).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*which is not ideal*\") {\n        $target = $p\n        break\n    }\n}\n\n$r = $target.Range\n$r.MoveEnd(1, -1)          # exclude the trailing paragraph mark\n$langId = $r.LanguageID    # formatting (en-AU) to apply to the new run\n$r.Collapse(0)             # collapse to the end of the paragraph's text\n$insPos = $r.Start\n\n# Append the missing \".\" plus a one-character sentinel. Inserting just \".\"\n# would put the paragraph's end-of-text position exactly on the paragraph\n# mark boundary, and adding a bookmark collapsed at that exact boundary\n# isn't reliable; the sentinel keeps the insertion point comfortably inside\n# the paragraph while we create the bookmark, and is removed afterwards.\n$r.InsertAfter(\".~\")\n$r.LanguageID = $langId\n\n# 3) Re-insert the \"_GoBack\" bookmark right after the period we just added.\n$bmPos = $insPos + 1\n$bmRange = $d.Range($bmPos, $bmPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# 4) Remove the temporary sentinel character.\n$d.Range($bmPos, $bmPos + 1).Delete()\n"}
